# Fixed traj read real robot data
# Columns B, D, F are negated (sign flip); columns C, E are replaced with (pi/2 - old value).
# Applied directly as exact literal values to avoid any floating point drift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.9903733673316458
$ws.Range("C1").Value = -0.0002537796048658035
$ws.Range("D1").Value = -1.40006462893066
$ws.Range("E1").Value = 0.1709853632432805
$ws.Range("F1").Value = 1.570796390562869

$ws.Range("B2").Value = 1.03574262259447
$ws.Range("C2").Value = -0.00009196697643037226
$ws.Range("D2").Value = -1.398862850774126
$ws.Range("E2").Value = 0.1720253258411236
$ws.Range("F2").Value = 1.57079639084348

$ws.Range("B3").Value = 1.239026965236186
$ws.Range("C3").Value = 0.0006330608780510351
$ws.Range("D3").Value = -1.393478087873858
$ws.Range("E3").Value = 0.1766850477574261
$ws.Range("F3").Value = 1.570796392100804

$ws.Range("B4").Value = 1.525028734600179
$ws.Range("C4").Value = 0.001653106231889857
$ws.Range("D4").Value = -1.385902237726871
$ws.Range("E4").Value = 0.1832408340786359
$ws.Range("F4").Value = 1.57079639386974

$ws.Range("B5").Value = 1.728313077241895
$ws.Range("C5").Value = 0.002378134086371267
$ws.Range("D5").Value = -1.380517474826603
$ws.Range("E5").Value = 0.1879005559949384
$ws.Range("F5").Value = 1.570796395127064

$ws.Range("B6").Value = 1.77368233250472
$ws.Range("C6").Value = 0.002539946714806697
$ws.Range("D6").Value = -1.379315696670069
$ws.Range("E6").Value = 0.1889405185927815
$ws.Range("F6").Value = 1.570796395407675
